# Insert a new weekly price-report row for Papa (Asterix, Provincia de
# Melipilla) into the middle of the "Terminal La Palmera de La Serena"
# Papa sheet, pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 497; everything currently at row 497
# and below shifts down to 498 and below (dimension grows to R599).
$ws.Rows("497:497").Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A497").Value = 8
$ws.Range("B497").Value = "Terminal La Palmera de La Serena"
$ws.Range("C497").Value = "Coquimbo"
$ws.Range("D497").Value = 44995
$ws.Range("E497").Value = 4
$ws.Range("F497").Value = 100114001
$ws.Range("G497").Value = "Papa"
$ws.Range("H497").Value = "Asterix"
$ws.Range("I497").Value = "1a (cosecha)"
$ws.Range("J497").Value = 1800
$ws.Range("K497").Value = 11000
$ws.Range("L497").Value = 12000
$ws.Range("M497").Value = 11500
$ws.Range("N497").Value = "`$/saco 25 kilos"
$ws.Range("O497").Value = "Provincia de Melipilla"
$ws.Range("P497").Value = 460
$ws.Range("Q497").Value = 25
$ws.Range("R497").Value = "Hortaliza"
